# "Added Entry For Steven Carl"
#
# Appends a new icebreaker-discussion entry to the end of the document,
# after David Singletary's example entry, without touching any existing
# content:
#
#   (blank paragraph)
#   9/7/2023 Steven Carl Hello My Name is Steven Carl. I was born and
#   raised here in Jacksonville Florida.

$d = $word.ActiveDocument

# Collapse a range to the very end of the document's main story.
$endRange = $d.Content
$endRange.Collapse(0)   # wdCollapseEnd

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# A blank paragraph, then the new entry split across two runs (the first
# run carries the date/name in Calibri, the second the free-form intro
# text typed with whatever the insertion-point default font was).
$newContentXml =
  "<w:p $wordNs>" +
    "<w:r><w:rPr><w:rFonts w:ascii=`"Calibri`"/></w:rPr></w:r>" +
  "</w:p>" +
  "<w:p $wordNs>" +
    "<w:r>" +
      "<w:rPr><w:rFonts w:ascii=`"Calibri`"/></w:rPr>" +
      "<w:t xml:space=`"preserve`">9/7/2023 Steven Carl </w:t>" +
    "</w:r>" +
    "<w:r>" +
      "<w:t xml:space=`"preserve`">Hello My Name is Steven Carl. I was born and raised here in Jacksonville Florida. </w:t>" +
    "</w:r>" +
  "</w:p>"

$endRange.InsertXML($newContentXml) | Out-Null
